$d = $word.ActiveDocument

# --- Locate the insertion point -----------------------------------------
# The new content is inserted right before the (empty) paragraph that sits
# 3 paragraphs ahead of the paragraph containing "5/22/20" - i.e. right
# after the second blank paragraph following the title line. We find the
# anchor text with Find (robust to any paraId/rsid noise) and then walk
# the Paragraphs collection to get a stable index instead of hard-coding
# paragraph numbers.
$find = $d.Content
$found = $find.Find.Execute("5/22/20", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find anchor text '5/22/20'"
}

$paras = $d.Paragraphs
$anchorIndex = -1
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Start -le $find.Start -and $find.Start -lt $p.Range.End) {
        $anchorIndex = $i
        break
    }
}
if ($anchorIndex -eq -1) {
    throw "Could not resolve anchor paragraph index"
}

$insertIndex = $anchorIndex - 3
$targetPara = $paras.Item($insertIndex)

# Collapse to a zero-length point right before the target paragraph's own
# content (rather than using the whole paragraph Range, which includes its
# end-of-paragraph mark and would get clobbered/replaced by InsertXML).
$targetRange = $targetPara.Range.Duplicate
$targetRange.Collapse(1)   # wdCollapseStart

# --- Build the new paragraphs as WordprocessingML and insert them -------
# Using InsertXML (rather than typing text) lets us reproduce the exact
# run-splitting and <w:proofErr/> spell-check markers from the source
# edit, including the curly quotes used verbatim in the notes.
$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$p1 = "<w:p $ns/>"
$p2 = "<w:p $ns><w:r><w:t>6/29/20</w:t></w:r></w:p>"
$p3 = "<w:p $ns/>"
$p4 = "<w:p $ns/>"
$p5 = "<w:p $ns>" +
      '<w:r><w:t xml:space="preserve">For IO behavior, added feature of ‘cutfirst3trials’ so we have files that include first three SG and those that automatically exclude (I think auto exclude is what we had before), note that I think for </w:t></w:r>' +
      '<w:proofErr w:type="spellStart"/>' +
      '<w:r><w:t>ephys</w:t></w:r>' +
      '<w:proofErr w:type="spellEnd"/>' +
      '<w:r><w:t xml:space="preserve"> we want to continue to exclude things but for behavioral analysis itself we can be more flexible.</w:t></w:r>' +
      '</w:p>'
$p6 = "<w:p $ns/>"
$p7 = "<w:p $ns/>"
$p8 = "<w:p $ns><w:r><w:t>6/24/20</w:t></w:r></w:p>"
$p9 = "<w:p $ns><w:r><w:t>Going to need to re-run initial processing code because I want to check out some additional analyses.</w:t></w:r></w:p>"
$p10 = "<w:p $ns>" +
       '<w:r><w:t>Wrote ‘behavior_V2’</w:t></w:r>' +
       '<w:r><w:t xml:space="preserve"> and V2 for </w:t></w:r>' +
       '<w:proofErr w:type="spellStart"/>' +
       '<w:r><w:t>io</w:t></w:r>' +
       '<w:proofErr w:type="spellEnd"/>' +
       '<w:r><w:t xml:space="preserve"> analysis</w:t></w:r>' +
       '</w:p>'

# Inserting a block of paragraphs at a collapsed point makes the *last*
# inserted paragraph mark merge with / take over the paragraph that used
# to start at that point (matching how Word's own paste/insert-XML works)
# - i.e. without padding, the existing (empty) paragraph that originally
# sat at $targetRange would be silently swallowed into our last new
# paragraph. Append one extra empty paragraph so that merge consumes our
# own padding instead of eating the paragraph that must survive unchanged
# right after our inserted block.
$padding = "<w:p $ns/>"

$xml = $p1 + $p2 + $p3 + $p4 + $p5 + $p6 + $p7 + $p8 + $p9 + $p10 + $padding

$targetRange.InsertXML($xml)
